# "change filename + fenetre saccade"
# Updates the saccade-detection windows on the first two image sheets
# (img_1.jpg / img_2.jpg): sheet1 gains one extra detection row (row 12)
# and several existing rows are rewritten with new X/Y/time values;
# sheet2's rows are rewritten in place (values shift by one "window"),
# and row 18 becomes a "real" data row (copied formatting from row 17,
# replacing the old duplicate tail row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (img_1.jpg)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 12 is brand new — clone formatting from row 11 (same borders/font
# as the rest of the data rows) before writing its values.
$ws1.Range("A11:I11").Copy($ws1.Range("A12:I12"))

$ws1.Range("C5").Value = -10
$ws1.Range("D5").Value = 80
$ws1.Range("E5").Value = 1
$ws1.Range("F5").Value = 58
$ws1.Range("G5").Value = 66.20758056640625
$ws1.Range("H5").Value = 66.28766632080078
$ws1.Range("I5").Value = 0.08008150011301041

$ws1.Range("C6").Value = -4
$ws1.Range("D6").Value = 57
$ws1.Range("E6").Value = 110
$ws1.Range("F6").Value = 77
$ws1.Range("G6").Value = 66.28358459472656
$ws1.Range("H6").Value = 66.36567687988281
$ws1.Range("I6").Value = 0.08208800107240677

$ws1.Range("C7").Value = -3
$ws1.Range("D7").Value = 78
$ws1.Range("E7").Value = 5
$ws1.Range("F7").Value = 69
$ws1.Range("G7").Value = 66.757568359375
$ws1.Range("H7").Value = 66.79957580566406
$ws1.Range("I7").Value = 0.04201050102710724

$ws1.Range("C8").Value = 125
$ws1.Range("D8").Value = 89
$ws1.Range("E8").Value = 110
$ws1.Range("F8").Value = 84
$ws1.Range("G8").Value = 67.51958465576172
$ws1.Range("H8").Value = 67.56157684326172
$ws1.Range("I8").Value = 0.04199250042438507

$ws1.Range("C9").Value = -8
$ws1.Range("D9").Value = -23
$ws1.Range("E9").Value = -4
$ws1.Range("F9").Value = -89
$ws1.Range("G9").Value = 67.91165924072266
$ws1.Range("H9").Value = 67.95758056640625
$ws1.Range("I9").Value = 0.04591749981045723

$ws1.Range("C10").Value = -116
$ws1.Range("D10").Value = -106
$ws1.Range("E10").Value = -119
$ws1.Range("F10").Value = -91
$ws1.Range("G10").Value = 68.45174407958984
$ws1.Range("H10").Value = 68.49562835693359
$ws1.Range("I10").Value = 0.04388250038027763

$ws1.Range("C11").Value = 122
$ws1.Range("D11").Value = 110
$ws1.Range("E11").Value = 48
$ws1.Range("F11").Value = 40
$ws1.Range("G11").Value = 68.96994018554688
$ws1.Range("H11").Value = 69.01158142089844
$ws1.Range("I11").Value = 0.0416405014693737

$ws1.Range("A12").Value = 10
$ws1.Range("B12").Value = "Normal"
$ws1.Range("C12").Value = -3
$ws1.Range("D12").Value = 3
$ws1.Range("E12").Value = 67
$ws1.Range("F12").Value = -34
$ws1.Range("G12").Value = 70.48759460449219
$ws1.Range("H12").Value = 70.53157043457031
$ws1.Range("I12").Value = 0.04397099837660789

# ---------------------------------------------------------------------
# Sheet 2 (img_2.jpg)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 18 stops being a duplicate of row 14/19 and becomes the row that
# used to be row 19's predecessor data — clone row 17's formatting
# (border+font) onto row 18 first, matching the rest of the data rows.
$ws2.Range("A17").Copy($ws2.Range("A18"))

$ws2.Range("E3").Value = 54
$ws2.Range("F3").Value = -91
$ws2.Range("H3").Value = 53.35359191894531
$ws2.Range("I3").Value = 0.08204949647188187

$ws2.Range("C4").Value = 56
$ws2.Range("D4").Value = -91
$ws2.Range("E4").Value = -91
$ws2.Range("F4").Value = -43
$ws2.Range("G4").Value = 53.34962844848633
$ws2.Range("H4").Value = 53.431640625
$ws2.Range("I4").Value = 0.08201199769973755

$ws2.Range("C5").Value = 9
$ws2.Range("D5").Value = -100
$ws2.Range("E5").Value = -14
$ws2.Range("F5").Value = 104
$ws2.Range("G5").Value = 54.01604461669922
$ws2.Range("H5").Value = 54.05765151977539
$ws2.Range("I5").Value = 0.04160749912261963

$ws2.Range("D6").Value = -55
$ws2.Range("E6").Value = -67
$ws2.Range("F6").Value = -127
$ws2.Range("G6").Value = 54.81957626342773
$ws2.Range("H6").Value = 54.861572265625
$ws2.Range("I6").Value = 0.04199599847197533

$ws2.Range("C7").Value = -34
$ws2.Range("D7").Value = 51
$ws2.Range("E7").Value = -5
$ws2.Range("F7").Value = 124
$ws2.Range("G7").Value = 55.94355010986328
$ws2.Range("H7").Value = 55.98361587524414
$ws2.Range("I7").Value = 0.04006800055503845

$ws2.Range("B8").Value = "Normal"
$ws2.Range("C8").Value = 90
$ws2.Range("D8").Value = -31
$ws2.Range("E8").Value = 52
$ws2.Range("F8").Value = 41
$ws2.Range("G8").Value = 56.76163864135742
$ws2.Range("H8").Value = 56.80354309082031
$ws2.Range("I8").Value = 0.0419050008058548

$ws2.Range("C9").Value = 85
$ws2.Range("D9").Value = 107
$ws2.Range("E9").Value = 84
$ws2.Range("F9").Value = 109
$ws2.Range("G9").Value = 57.29579162597656
$ws2.Range("H9").Value = 57.29965972900391
$ws2.Range("I9").Value = 0.00386850000359118

$ws2.Range("C10").Value = 74
$ws2.Range("D10").Value = 110
$ws2.Range("E10").Value = 76
$ws2.Range("G10").Value = 57.30368804931641
$ws2.Range("H10").Value = 57.30768966674805
$ws2.Range("I10").Value = 0.004002499859780073

$ws2.Range("D11").Value = 109
$ws2.Range("E11").Value = 80
$ws2.Range("F11").Value = 107
$ws2.Range("G11").Value = 57.31179809570312
$ws2.Range("H11").Value = 57.3157844543457
$ws2.Range("I11").Value = 0.00398549996316433

$ws2.Range("C12").Value = 80
$ws2.Range("D12").Value = 106
$ws2.Range("E12").Value = 81
$ws2.Range("F12").Value = 103
$ws2.Range("G12").Value = 57.31963729858398
$ws2.Range("H12").Value = 57.32362365722656
$ws2.Range("I12").Value = 0.003988500218838453

$ws2.Range("B13").Value = "Micro"
$ws2.Range("C13").Value = 85
$ws2.Range("D13").Value = 107
$ws2.Range("E13").Value = 84
$ws2.Range("F13").Value = 104
$ws2.Range("G13").Value = 57.32758331298828
$ws2.Range("H13").Value = 57.33161163330078
$ws2.Range("I13").Value = 0.004026500042527914

$ws2.Range("C14").Value = 87
$ws2.Range("D14").Value = 104
$ws2.Range("E14").Value = 88
$ws2.Range("F14").Value = 108
$ws2.Range("G14").Value = 57.33960342407227
$ws2.Range("H14").Value = 57.38351058959961
$ws2.Range("I14").Value = 0.04390550032258034

$ws2.Range("C15").Value = 88
$ws2.Range("D15").Value = 112
$ws2.Range("E15").Value = 71
$ws2.Range("F15").Value = -48
$ws2.Range("G15").Value = 57.37959671020508
$ws2.Range("H15").Value = 57.46163177490234
$ws2.Range("I15").Value = 0.08203549683094025

$ws2.Range("C16").Value = 123
$ws2.Range("D16").Value = 119
$ws2.Range("E16").Value = -124
$ws2.Range("F16").Value = 92
$ws2.Range("G16").Value = 58.08774566650391
$ws2.Range("H16").Value = 58.13164138793945
$ws2.Range("I16").Value = 0.04389650002121925

$ws2.Range("C17").Value = 42
$ws2.Range("D17").Value = 36
$ws2.Range("E17").Value = 52
$ws2.Range("F17").Value = 38
$ws2.Range("G17").Value = 58.5036506652832
$ws2.Range("H17").Value = 58.54961395263672
$ws2.Range("I17").Value = 0.04596599936485291

$ws2.Range("A18").Value = 16
$ws2.Range("B18").Value = "Normal"
$ws2.Range("C18").Value = -128
$ws2.Range("D18").Value = 83
$ws2.Range("E18").Value = 11
$ws2.Range("F18").Value = -10
$ws2.Range("G18").Value = 59.18439102172852
$ws2.Range("H18").Value = 59.22768020629883
$ws2.Range("I18").Value = 0.04328899830579758

Write-Host "edit applied"
